# edit.ps1
# Applies the "Little_Medium" typography addition and board/debug-mode
# indicator translation rows to the TouchGFX texts workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Typography sheet: add the new "Little_Medium" typography entry
#    in row 11 (the first free row of the Typography table).
# ---------------------------------------------------------------
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("B11").Value = "Little_Medium"
$wsTypo.Range("C11").Value = "Roboto-Regular.ttf"
$wsTypo.Range("D11").Value = 35
$wsTypo.Range("E11").Value = 4
$wsTypo.Range("F11").Value = "?"
$wsTypo.Range("H11").Value = "42-91"

# ---------------------------------------------------------------
# 2) Translation sheet: rename every existing use of the
#    "LittleMedium" typography to the newly introduced
#    "Little_Medium" typography.
# ---------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Columns.Item(3).Replace("LittleMedium", "Little_Medium", 1, 1, $false, $false, $false, $false)

# ---------------------------------------------------------------
# 3) Translation sheet: populate the previously empty rows
#    (266-295) with the board debug-mode indicator texts that now
#    use the "Little_Medium" typography.
# ---------------------------------------------------------------
$rows = @(
    @{Row=266; B="SingleUseId286"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=267; B="SingleUseId287"; C="Little_Medium"; D="Left"; E="ind_6"; F="LTR"},
    @{Row=268; B="SingleUseId288"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=269; B="SingleUseId289"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=270; B="SingleUseId290"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=271; B="SingleUseId291"; C="Little_Medium"; D="Left"; E="ind_7"; F="LTR"},
    @{Row=272; B="SingleUseId292"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=273; B="SingleUseId293"; C="Little_Medium"; D="Left"; E="ind_6`n"; F="LTR"},
    @{Row=274; B="SingleUseId294"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=275; B="SingleUseId295"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=276; B="SingleUseId296"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=277; B="SingleUseId297"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=278; B="SingleUseId298"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=279; B="SingleUseId299"; C="Little_Medium"; D="Left"; E="ind_7"; F="LTR"},
    @{Row=280; B="SingleUseId300"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=281; B="SingleUseId301"; C="Little_Medium"; D="Left"; E="ind_6"; F="LTR"},
    @{Row=282; B="SingleUseId302"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=283; B="SingleUseId303"; C="Little_Medium"; D="Left"; E="ind_8`n"; F="LTR"},
    @{Row=284; B="SingleUseId304"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=285; B="SingleUseId305"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=286; B="SingleUseId306"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=287; B="SingleUseId307"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=288; B="SingleUseId308"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=289; B="SingleUseId309"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=290; B="SingleUseId310"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=291; B="SingleUseId311"; C="Little_Medium"; D="Left"; E="ind_8`n"; F="LTR"},
    @{Row=292; B="SingleUseId312"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=293; B="SingleUseId313"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"},
    @{Row=294; B="SingleUseId314"; C="Little_Medium"; D="Center"; E="<value>"; F="LTR"},
    @{Row=295; B="SingleUseId315"; C="Little_Medium"; D="Left"; E="0.00"; F="LTR"}
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $wsTrans.Range("B$rowNum").Value = $r.B
    $wsTrans.Range("C$rowNum").Value = $r.C
    $wsTrans.Range("D$rowNum").Value = $r.D

    # The "GB" column sometimes holds number-looking labels (e.g. "0.00")
    # that must stay plain text, exactly like the rest of the table, so
    # force text formatting only for those to avoid Excel auto-converting
    # them to real numbers.
    $eCell = $wsTrans.Range("E$rowNum")
    if ($r.E -match '^-?[0-9]+(\.[0-9]+)?$') {
        $eCell.NumberFormat = "@"
    }
    $eCell.Value = $r.E

    $wsTrans.Range("F$rowNum").Value = $r.F
}
